$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 30000
$ws.Range("J3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30228
# Row 64
$ws.Range("H64").Value = 9999
$ws.Range("I64").Value = 9999.143
$ws.Range("J64").Value = 9998
$ws.Range("K64").Value = 9999.143
$ws.Range("L64").Value = 9998
$ws.Range("M64").Value = -9751.143
$ws.Range("N64").Value = -10494
# Row 67
$ws.Range("H67").Value = 9999
$ws.Range("I67").Value = 9999.143
$ws.Range("J67").Value = 9998
$ws.Range("K67").Value = 9999.143
$ws.Range("L67").Value = 9998
$ws.Range("M67").Value = -9141.143
$ws.Range("N67").Value = -11714
# Row 102
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -36490
# Row 132
$ws.Range("H132").Value = 2553.4517
$ws.Range("I132").Value = 2538.5667
$ws.Range("K132").Value = 7615.7001
$ws.Range("M132").Value = -5085.7001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 474
$ws.Range("I97").Value = 477.625
$ws.Range("K97").Value = 477.625
$ws.Range("M97").Value = 18.375
# Row 122
$ws.Range("H122").Value = 3449.647
$ws.Range("I122").Value = 3218.9614
$ws.Range("J122").Value = 4199.375
$ws.Range("K122").Value = 9656.8842
$ws.Range("L122").Value = 12598.125
$ws.Range("M122").Value = -7206.8842
$ws.Range("N122").Value = -17498.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1525.25
$ws.Range("I99").Value = 1386.2142
$ws.Range("J99").Value = 2498.5
$ws.Range("K99").Value = 1386.2142
$ws.Range("L99").Value = 2498.5
$ws.Range("M99").Value = 111.7858000000001
$ws.Range("N99").Value = -5494.5
# Row 107
$ws.Range("H107").Value = 5313.722
$ws.Range("I107").Value = 5313.722
$ws.Range("K107").Value = 5313.722
$ws.Range("M107").Value = -3393.722
# Row 134
$ws.Range("H134").Value = 2547.923
$ws.Range("I134").Value = 1749.2609
$ws.Range("K134").Value = 5247.7827
$ws.Range("M134").Value = -2712.7827

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 11021.1875
$ws.Range("I31").Value = 3180
$ws.Range("K31").Value = 3180
$ws.Range("M31").Value = -2885
# Row 34
$ws.Range("H34").Value = 11021.1875
$ws.Range("I34").Value = 3180
$ws.Range("K34").Value = 3180
$ws.Range("M34").Value = -2978
# Row 105
$ws.Range("H105").Value = 3174.5
$ws.Range("I105").Value = 3566
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 3566
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -1819
$ws.Range("N105").Value = -5494
# Row 107
$ws.Range("H107").Value = 2230.1667
$ws.Range("I107").Value = 1984.7778
$ws.Range("J107").Value = 2966.3333
$ws.Range("K107").Value = 1984.7778
$ws.Range("L107").Value = 2966.3333
$ws.Range("M107").Value = -64.77780000000007
$ws.Range("N107").Value = -6806.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 110
$ws.Range("H110").Value = 4197.8335
$ws.Range("I110").Value = 4197.8335
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 12593.5005
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -8503.500499999998
$ws.Range("N110").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 53839
$ws.Range("J26").Value = 53839
$ws.Range("L26").Value = 53839
$ws.Range("N26").Value = -54399
# Row 50
$ws.Range("H50").Value = 53839
$ws.Range("J50").Value = 53839
$ws.Range("L50").Value = 53839
$ws.Range("N50").Value = -54835
# Row 70
$ws.Range("H70").Value = 7211.0303
$ws.Range("I70").Value = 5997.5
$ws.Range("K70").Value = 5997.5
$ws.Range("M70").Value = -5727.5
# Row 73
$ws.Range("H73").Value = 7211.0303
$ws.Range("I73").Value = 5997.5
$ws.Range("K73").Value = 5997.5
$ws.Range("M73").Value = -5061.5
# Row 102
$ws.Range("H102").Value = 2586.2693
$ws.Range("I102").Value = 1655.1818
$ws.Range("J102").Value = 7707.25
$ws.Range("K102").Value = 1655.1818
$ws.Range("L102").Value = 7707.25
$ws.Range("M102").Value = -33.18180000000007
$ws.Range("N102").Value = -10951.25
# Row 122
$ws.Range("H122").Value = 3189.72
$ws.Range("J122").Value = 6414.6665
$ws.Range("L122").Value = 19243.9995
$ws.Range("N122").Value = -24143.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1322.5714
$ws.Range("I22").Value = 1174.0834
$ws.Range("J22").Value = 1520.5555
$ws.Range("K22").Value = 1174.0834
$ws.Range("L22").Value = 1520.5555
$ws.Range("M22").Value = -879.0834
$ws.Range("N22").Value = -2110.5555
# Row 27
$ws.Range("H27").Value = 1322.5714
$ws.Range("I27").Value = 1174.0834
$ws.Range("J27").Value = 1520.5555
$ws.Range("K27").Value = 1174.0834
$ws.Range("L27").Value = 1520.5555
$ws.Range("M27").Value = -1067.0834
$ws.Range("N27").Value = -1734.5555
# Row 68
$ws.Range("H68").Value = 4541.143
$ws.Range("I68").Value = 4580
$ws.Range("K68").Value = 4580
$ws.Range("M68").Value = -3831
# Row 71
$ws.Range("H71").Value = 4541.143
$ws.Range("I71").Value = 4580
$ws.Range("K71").Value = 22901
$ws.Range("M71").Value = -19156
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
# Row 100
$ws.Range("H100").Value = 2335
$ws.Range("I100").Value = 1002.5
$ws.Range("K100").Value = 1002.5
$ws.Range("M100").Value = -461.5
# Row 132
$ws.Range("H132").Value = 3936.5293
$ws.Range("I132").Value = 2274.7273
$ws.Range("K132").Value = 6824.1819
$ws.Range("M132").Value = -4294.1819
# Row 136
$ws.Range("H136").Value = 5261.956
$ws.Range("I136").Value = 2626.853
$ws.Range("K136").Value = 7880.559
$ws.Range("M136").Value = -5330.559

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 42
$ws.Range("H42").Value = 64500
$ws.Range("J42").Value = 76000
$ws.Range("L42").Value = 76000
$ws.Range("N42").Value = -76756
# Row 62
$ws.Range("H62").Value = 15250
$ws.Range("J62").Value = 15250
$ws.Range("L62").Value = 15250
$ws.Range("N62").Value = -16498
# Row 65
$ws.Range("H65").Value = 15250
$ws.Range("J65").Value = 15250
$ws.Range("L65").Value = 76250
$ws.Range("N65").Value = -82490
# Row 81
$ws.Range("H81").Value = 1349.4445
$ws.Range("I81").Value = 1535.6364
$ws.Range("J81").Value = 1056.8572
$ws.Range("K81").Value = 3071.2728
$ws.Range("L81").Value = 2113.7144
$ws.Range("M81").Value = -2010.2728
$ws.Range("N81").Value = -4235.7144
# Row 84
$ws.Range("H84").Value = 1349.4445
$ws.Range("I84").Value = 1535.6364
$ws.Range("J84").Value = 1056.8572
$ws.Range("K84").Value = 15356.364
$ws.Range("L84").Value = 10568.572
$ws.Range("M84").Value = -10052.364
$ws.Range("N84").Value = -21176.572
# Row 96
$ws.Range("H96").Value = 5153.6665
$ws.Range("I96").Value = 2691.2856
$ws.Range("J96").Value = 6384.857
$ws.Range("K96").Value = 2691.2856
$ws.Range("L96").Value = 6384.857
$ws.Range("M96").Value = -1318.2856
$ws.Range("N96").Value = -9130.857
# Row 122
$ws.Range("H122").Value = 1789.7142
$ws.Range("I122").Value = 1791.6818
$ws.Range("J122").Value = 1772.4
$ws.Range("K122").Value = 5375.0454
$ws.Range("L122").Value = 5317.200000000001
$ws.Range("M122").Value = -2925.0454
$ws.Range("N122").Value = -10217.2
# Row 126
$ws.Range("H126").Value = 1214.6207
$ws.Range("I126").Value = 1175.66
$ws.Range("J126").Value = 1458.125
$ws.Range("K126").Value = 3526.98
$ws.Range("L126").Value = 4374.375
$ws.Range("M126").Value = -1056.98
$ws.Range("N126").Value = -9314.375
# Row 132
$ws.Range("H132").Value = 2878.805
$ws.Range("I132").Value = 1932.9714
$ws.Range("K132").Value = 5798.914199999999
$ws.Range("M132").Value = -3268.914199999999
# Row 136
$ws.Range("H136").Value = 2862.8125
$ws.Range("I136").Value = 671.8889
$ws.Range("K136").Value = 2015.6667
$ws.Range("M136").Value = 534.3332999999998

